$wb = $excel.ActiveWorkbook

# --- "login" sheet -------------------------------------------------
# The test-data rows were bumped: Palatine8/9 -> Palatine10/11,
# and (on the "order" sheet) ChicagoCity7/8 -> ChicagoCity9/10.
$wsLogin = $wb.Worksheets.Item("login")

$wsLogin.Range("F2").Value = "AutomationTestPalatine10"
$wsLogin.Range("H2").Value = "AutomationTestPalatine10@gmail.com"
$wsLogin.Range("F3").Value = "AutomationTestPalatine11"
$wsLogin.Range("H3").Value = "AutomationTestPalatine11@gmail.com"

# --- "order" sheet ---------------------------------------------------
$wsOrder = $wb.Worksheets.Item("order")

$wsOrder.Range("R2").Value = "TestChicagoCity9"
$wsOrder.Range("S2").Value = "TestChicagoCity9"
$wsOrder.Range("T2").Value = "TestChicagoCity9@gmail.com"
$wsOrder.Range("R3").Value = "TestChicagoCity10"
$wsOrder.Range("S3").Value = "TestChicagoCity10"
$wsOrder.Range("T3").Value = "TestChicagoCity10@gmail.com"

# --- sheet view / selection state ------------------------------------
# login!K4 selected (scrolled so column F is left-most visible),
# order!U9 selected; order stays the active tab, matching the diff.
$wsLogin.Range("K4").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1

$wsOrder.Range("U9").Select()
